# New Member Beneficiary Role scripts
# Adds 5 new test-case rows to the "General" sheet and 5 corresponding
# detail rows (with beneficiary/loan/CC/MM/CD/personal-loan role data)
# to the "DataTwo" sheet, then leaves the workbook with DataTwo active.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "DataTwo" sheet (sheet5) - new detail rows 7..11
# ---------------------------------------------------------------------
$dataTwo = $wb.Worksheets.Item("DataTwo")

# Rows 7, 9 and 10 only use columns A..S (same shape as row 6's template).
$dataTwo.Range("A6:AE6").Copy($dataTwo.Range("A7:AE7"))
$dataTwo.Range("A6:AE6").Copy($dataTwo.Range("A9:AE9"))
$dataTwo.Range("A6:AE6").Copy($dataTwo.Range("A10:AE10"))

# Rows 8 and 11 also populate T..AA, so use row 2's template (it already
# has those columns formatted).
$dataTwo.Range("A2:AA2").Copy($dataTwo.Range("A8:AA8"))
$dataTwo.Range("A2:AA2").Copy($dataTwo.Range("A11:AA11"))

# --- Row 7: C24336_NonMemberBeneficiarySavingsRole ---
$dataTwo.Range("A7").Value = "C24336_NonMemberBeneficiarySavingsRole"
$dataTwo.Range("B7").Value = "Peter"
$dataTwo.Range("C7").Value = "Ford"
$dataTwo.Range("F7").Value = "'46889712"
$dataTwo.Range("G7").Value = "'06/06/1955"
$dataTwo.Range("R7").Value = "test+6@email.com"
$dataTwo.Range("S7").Value = "512365541"

# --- Row 8: C24335_NonMemberJOLoanandCCRoles ---
$dataTwo.Range("A8").Value = "C24335_NonMemberJOLoanandCCRoles"
$dataTwo.Range("B8").Value = 500
$dataTwo.Range("C8").Value = 4000
$dataTwo.Range("F8").Value = "'52422589"
$dataTwo.Range("G8").Value = "'07/07/1965"
$dataTwo.Range("R8").Value = "test+7@email.com"
$dataTwo.Range("S8").Value = "2000"
$dataTwo.Range("T8").Value = "3"
$dataTwo.Range("U8").Value = "5000"
$dataTwo.Range("V8").Value = "Barb"
$dataTwo.Range("W8").Value = "Henson"
$dataTwo.Range("X8").Value = "522365541"
$dataTwo.Range("Y8").Value = "Classic MasterCard"
$dataTwo.Range("Z8").Value = "My Way Loan: $4,000"
$dataTwo.Range("AA8").Value = "Share"

# --- Row 9: C24334_NonMemberBeneficiaryMMRole ---
$dataTwo.Range("A9").Value = "C24334_NonMemberBeneficiaryMMRole"
$dataTwo.Range("B9").Value = "Julie"
$dataTwo.Range("C9").Value = "Morris"
$dataTwo.Range("F9").Value = "'48889712"
$dataTwo.Range("G9").Value = "'08/08/1975"
$dataTwo.Range("R9").Value = "test+8@email.com"
$dataTwo.Range("S9").Value = "532365541"

# --- Row 10: C24333_NonMemberBeneficiaryCDRole ---
$dataTwo.Range("A10").Value = "C24333_NonMemberBeneficiaryCDRole"
$dataTwo.Range("B10").Value = "Todd"
$dataTwo.Range("C10").Value = "Lambert"
$dataTwo.Range("F10").Value = "'49889712"
$dataTwo.Range("G10").Value = "'09/09/1985"
$dataTwo.Range("R10").Value = "test+9@email.com"
$dataTwo.Range("S10").Value = "542365541"

# --- Row 11: C24332_NonMemberBeneficiaryPersonalLoanRole ---
$dataTwo.Range("A11").Value = "C24332_NonMemberBeneficiaryPersonalLoanRole"
$dataTwo.Range("B11").Value = "Kelly"
$dataTwo.Range("C11").Value = "Colt"
$dataTwo.Range("F11").Value = "'53422589"
$dataTwo.Range("G11").Value = "'09/09/1985"
$dataTwo.Range("R11").Value = "test@email.com"
$dataTwo.Range("S11").Value = "2000"
$dataTwo.Range("T11").Value = "3"
$dataTwo.Range("U11").Value = "5000"
$dataTwo.Range("V11").Value = 1000
$dataTwo.Range("W11").Value = "552365541"

# ---------------------------------------------------------------------
# 2) "General" sheet (sheet1) - new TestCaseID rows 95..99
# ---------------------------------------------------------------------
$general = $wb.Worksheets.Item("General")

# Rows 95..98 already exist as (mostly blank) template rows; row 99 is
# brand new, so seed it from row 94's layout first.
$general.Range("A94:D94").Copy($general.Range("A99:D99"))

$general.Range("A95").Value = "C24336_NonMemberBeneficiarySavingsRole"
$general.Range("B95").Value = "https://forms-preprod.fivision.com/tdecu/oa/Default.aspx"
$general.Range("C95").Value = "Yes"
$general.Range("D95").Value = "Chrome"

$general.Range("A96").Value = "C24335_NonMemberJOLoanandCCRoles"
$general.Range("B96").Value = "https://forms-preprod.fivision.com/tdecu/oa/Default.aspx"
$general.Range("C96").Value = "Yes"
$general.Range("D96").Value = "Chrome"

$general.Range("A97").Value = "C24334_NonMemberBeneficiaryMMRole"
$general.Range("B97").Value = "https://forms-preprod.fivision.com/tdecu/oa/Default.aspx"
$general.Range("C97").Value = "Yes"
$general.Range("D97").Value = "Chrome"

$general.Range("A98").Value = "C24333_NonMemberBeneficiaryCDRole"
$general.Range("B98").Value = "https://forms-preprod.fivision.com/tdecu/oa/Default.aspx"
$general.Range("C98").Value = "Yes"
$general.Range("D98").Value = "Chrome"

$general.Range("A99").Value = "C24332_NonMemberBeneficiaryPersonalLoanRole"
$general.Range("B99").Value = "https://forms-preprod.fivision.com/tdecu/oa/Default.aspx"
$general.Range("C99").Value = "Yes"
$general.Range("D99").Value = "Chrome"

# Hyperlinks for the new rows' URL cells (B95..B99), mirroring every
# other row in this column.
$general.Hyperlinks.Add($general.Range("B95"), "https://forms-preprod.fivision.com/tdecu/oa/Default.aspx")
$general.Hyperlinks.Add($general.Range("B96"), "https://forms-preprod.fivision.com/tdecu/oa/Default.aspx")
$general.Hyperlinks.Add($general.Range("B97"), "https://forms-preprod.fivision.com/tdecu/oa/Default.aspx")
$general.Hyperlinks.Add($general.Range("B98"), "https://forms-preprod.fivision.com/tdecu/oa/Default.aspx")
$general.Hyperlinks.Add($general.Range("B99"), "https://forms-preprod.fivision.com/tdecu/oa/Default.aspx")

# ---------------------------------------------------------------------
# 3) View state - DataTwo becomes the active sheet/tab, with the new
#    last row selected on each touched sheet.
# ---------------------------------------------------------------------
$general.Range("B99:D99").Select()

$dataTwo.Activate()
$dataTwo.Range("A11").Select()
